$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values per diff
$ws.Range("L4").Value = 7
$ws.Range("L10").Value = 10

# Update sheet view: topLeftCell and selection
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("I11").Select()
